$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot (scheduled GitHub Actions update).
# "Price" values that look numeric (e.g. 156.60, 0.07034) are written with a
# leading apostrophe so Excel stores them as literal text -- preserving exact
# trailing zeros / dotted-thousands formatting -- instead of auto-converting
# them to numbers (which would drop precision, e.g. 156.60 -> 156.6).

$ws.Range("D2").Value = "28.940.90"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "1.904.82"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'325.03"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.4597"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").Value = "'0.3818"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").Value = "'45.57"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "'0.07735"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'0.9823"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "'22.08"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").Value = "1.944.21"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").Value = "'6.999"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "'5.686"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "'0.07034"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'84.26"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").Value = "'0.000009557"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'16.75"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "28.953.82"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "'5.334"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "'10.96"
$ws.Range("D25").Value = "2.162.39"
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("D26").Value = "'2.076"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").Value = "'156.60"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'19.19"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "'5.612"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("D30").Value = "'117.62"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -6.83%  "
$ws.Range("D32").Value = "'0.09263"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").Value = "'0.8635"
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").Value = "'5.113"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").Value = "'1.255"
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("D36").Value = "'3.018"
$ws.Range("E36").Value = "  -5.09%  "
$ws.Range("D37").Value = "'0.05719"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").Value = "'1.143"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "'1.003"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D41").Value = "'7.494"
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("D42").Value = "'0.5538"
$ws.Range("E42").Value = "  -3.88%  "
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").Value = "'9.324"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "'2.761"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").Value = "'0.5221"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  -6.33%  "
$ws.Range("D48").Value = "'2.107"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").Value = "'0.000002643"
$ws.Range("E49").Value = "  -20.00%  "
$ws.Range("D50").Value = "'0.06816"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").Value = "'112.11"
$ws.Range("E51").Value = "  -2.18%  "
